$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by exactly one day
for ($r = 2; $r -le 97; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value2 = $cur + 1
}

# Update the Notified Production (MW) values for rows 2-93 (rows 94-97 remain 0)
$ws.Cells.Item(2, 2).Value = 392
$ws.Cells.Item(3, 2).Value = 392
$ws.Cells.Item(4, 2).Value = 400
$ws.Cells.Item(5, 2).Value = 399
$ws.Cells.Item(6, 2).Value = 413
$ws.Cells.Item(7, 2).Value = 409
$ws.Cells.Item(8, 2).Value = 407
$ws.Cells.Item(9, 2).Value = 406
$ws.Cells.Item(10, 2).Value = 426
$ws.Cells.Item(11, 2).Value = 427
$ws.Cells.Item(12, 2).Value = 426
$ws.Cells.Item(13, 2).Value = 425
$ws.Cells.Item(14, 2).Value = 416
$ws.Cells.Item(15, 2).Value = 414
$ws.Cells.Item(16, 2).Value = 408
$ws.Cells.Item(17, 2).Value = 390
$ws.Cells.Item(18, 2).Value = 368
$ws.Cells.Item(19, 2).Value = 365
$ws.Cells.Item(20, 2).Value = 356
$ws.Cells.Item(21, 2).Value = 353
$ws.Cells.Item(22, 2).Value = 337
$ws.Cells.Item(23, 2).Value = 329
$ws.Cells.Item(24, 2).Value = 323
$ws.Cells.Item(25, 2).Value = 320
$ws.Cells.Item(26, 2).Value = 302
$ws.Cells.Item(27, 2).Value = 299
$ws.Cells.Item(28, 2).Value = 296
$ws.Cells.Item(29, 2).Value = 292
$ws.Cells.Item(30, 2).Value = 288
$ws.Cells.Item(31, 2).Value = 285
$ws.Cells.Item(32, 2).Value = 281
$ws.Cells.Item(33, 2).Value = 277
$ws.Cells.Item(34, 2).Value = 290
$ws.Cells.Item(35, 2).Value = 281
$ws.Cells.Item(36, 2).Value = 275
$ws.Cells.Item(37, 2).Value = 269
$ws.Cells.Item(38, 2).Value = 270
$ws.Cells.Item(39, 2).Value = 269
$ws.Cells.Item(40, 2).Value = 275
$ws.Cells.Item(41, 2).Value = 277
$ws.Cells.Item(42, 2).Value = 307
$ws.Cells.Item(43, 2).Value = 322
$ws.Cells.Item(44, 2).Value = 333
$ws.Cells.Item(45, 2).Value = 343
$ws.Cells.Item(46, 2).Value = 372
$ws.Cells.Item(47, 2).Value = 377
$ws.Cells.Item(48, 2).Value = 382
$ws.Cells.Item(49, 2).Value = 388
$ws.Cells.Item(50, 2).Value = 400
$ws.Cells.Item(51, 2).Value = 402
$ws.Cells.Item(52, 2).Value = 402
$ws.Cells.Item(53, 2).Value = 403
$ws.Cells.Item(54, 2).Value = 410
$ws.Cells.Item(55, 2).Value = 412
$ws.Cells.Item(56, 2).Value = 413
$ws.Cells.Item(57, 2).Value = 414
$ws.Cells.Item(58, 2).Value = 436
$ws.Cells.Item(59, 2).Value = 442
$ws.Cells.Item(60, 2).Value = 446
$ws.Cells.Item(61, 2).Value = 450
$ws.Cells.Item(62, 2).Value = 487
$ws.Cells.Item(63, 2).Value = 498
$ws.Cells.Item(64, 2).Value = 508
$ws.Cells.Item(65, 2).Value = 517
$ws.Cells.Item(66, 2).Value = 591
$ws.Cells.Item(67, 2).Value = 612
$ws.Cells.Item(68, 2).Value = 633
$ws.Cells.Item(69, 2).Value = 655
$ws.Cells.Item(70, 2).Value = 740
$ws.Cells.Item(71, 2).Value = 769
$ws.Cells.Item(72, 2).Value = 798
$ws.Cells.Item(73, 2).Value = 826
$ws.Cells.Item(74, 2).Value = 891
$ws.Cells.Item(75, 2).Value = 905
$ws.Cells.Item(76, 2).Value = 919
$ws.Cells.Item(77, 2).Value = 933
$ws.Cells.Item(78, 2).Value = 917
$ws.Cells.Item(79, 2).Value = 922
$ws.Cells.Item(80, 2).Value = 925
$ws.Cells.Item(81, 2).Value = 928
$ws.Cells.Item(82, 2).Value = 936
$ws.Cells.Item(83, 2).Value = 927
$ws.Cells.Item(84, 2).Value = 921
$ws.Cells.Item(85, 2).Value = 917
$ws.Cells.Item(86, 2).Value = 868
$ws.Cells.Item(87, 2).Value = 862
$ws.Cells.Item(88, 2).Value = 853
$ws.Cells.Item(89, 2).Value = 844
$ws.Cells.Item(90, 2).Value = 782
$ws.Cells.Item(91, 2).Value = 769
$ws.Cells.Item(92, 2).Value = 757
$ws.Cells.Item(93, 2).Value = 745
